# backwardElimination.xlsx edit
#
# The underlying statsmodels OLS "summary()" text blocks (stored verbatim in
# column B, row 2 of each of the 29 worksheets) were produced by re-running
# the backward-elimination script the next day. Only the "Date:" and "Time:"
# header lines inside that captured console text change; every number in the
# regression tables is identical to the previous run.
#
# Date goes from "Wed, 01 Jan 2020" -> "Thu, 02 Jan 2020" on every sheet.
# Time is sheet-specific (each step of the elimination ran a touch later).

$wb = $excel.ActiveWorkbook

$oldDate = "Wed, 01 Jan 2020"
$newDate = "Thu, 02 Jan 2020"

$timeMap = @(
    @{ Sheet = 1; OldTime = "23:18:40"; NewTime = "20:48:31" }
    @{ Sheet = 2; OldTime = "23:18:41"; NewTime = "20:48:31" }
    @{ Sheet = 3; OldTime = "23:18:41"; NewTime = "20:48:31" }
    @{ Sheet = 4; OldTime = "23:18:41"; NewTime = "20:48:31" }
    @{ Sheet = 5; OldTime = "23:18:41"; NewTime = "20:48:31" }
    @{ Sheet = 6; OldTime = "23:18:41"; NewTime = "20:48:31" }
    @{ Sheet = 7; OldTime = "23:18:41"; NewTime = "20:48:31" }
    @{ Sheet = 8; OldTime = "23:18:41"; NewTime = "20:48:31" }
    @{ Sheet = 9; OldTime = "23:18:41"; NewTime = "20:48:31" }
    @{ Sheet = 10; OldTime = "23:18:41"; NewTime = "20:48:31" }
    @{ Sheet = 11; OldTime = "23:18:41"; NewTime = "20:48:31" }
    @{ Sheet = 12; OldTime = "23:18:41"; NewTime = "20:48:31" }
    @{ Sheet = 13; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 14; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 15; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 16; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 17; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 18; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 19; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 20; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 21; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 22; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 23; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 24; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 25; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 26; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 27; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 28; OldTime = "23:18:41"; NewTime = "20:48:32" }
    @{ Sheet = 29; OldTime = "23:18:41"; NewTime = "20:48:32" }
)

foreach ($entry in $timeMap) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    $cell = $ws.Range("B2")

    $text = $cell.Value()

    if ($text.IndexOf($oldDate) -lt 0) {
        throw "Sheet $($entry.Sheet): expected date '$oldDate' not found"
    }
    if ($text.IndexOf($entry.OldTime) -lt 0) {
        throw "Sheet $($entry.Sheet): expected time '$($entry.OldTime)' not found"
    }

    $text = $text.Replace($oldDate, $newDate)
    $text = $text.Replace($entry.OldTime, $entry.NewTime)

    $cell.Value = $text
}
